# Automatic update of files.
# - Column C (rows 2-10): bump the "Förändrad" date serial by one day (46059 -> 46060)
# - Rows 4,5,7,8,9: the underlying data rows have shifted position; re-apply the
#   correct Beteckning/Datum/Area (ha) values for A, B and G columns so that each
#   row again shows the right record after the cyclic re-sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bump "Förändrad" (column C) by one day for every data row ---
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}

# --- Re-assign A (Beteckning), B (Datum) and G (Area (ha)) for the rows that moved ---
$ws.Cells.Item(4, 1).Value2 = "A 26262-2024"
$ws.Cells.Item(4, 2).Value2 = 45468.66077546297
$ws.Cells.Item(4, 7).Value2 = 0.6

$ws.Cells.Item(5, 1).Value2 = "A 14517-2023"
$ws.Cells.Item(5, 2).Value2 = 45012
$ws.Cells.Item(5, 7).Value2 = 0.6

$ws.Cells.Item(7, 1).Value2 = "A 14516-2023"
$ws.Cells.Item(7, 2).Value2 = 45012.86600694444
$ws.Cells.Item(7, 7).Value2 = 0.4

$ws.Cells.Item(8, 1).Value2 = "A 23798-2024"
$ws.Cells.Item(8, 2).Value2 = 45455.43208333333
$ws.Cells.Item(8, 7).Value2 = 1.3

$ws.Cells.Item(9, 1).Value2 = "A 4156-2023"
$ws.Cells.Item(9, 2).Value2 = 44953
$ws.Cells.Item(9, 7).Value2 = 1.5
